# #327 Ajout des profils d'acces
#
# 1. Bump the "Date" metadata value.
# 2. On the "Elements" sheet, swap the two mapping columns:
#    AK ("Mapping: RIM Mapping") <-> AL ("Mapping: Spécification métier
#    vers l'extension ROR HealthcareServicePatientType"), header included,
#    and swap their column widths to match.

$wb = $excel.ActiveWorkbook

# --- 1. Metadata!B8 : refresh the generation Date ---------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value2 = "2024-03-19T13:17:15+00:00"

# --- 2. Elements sheet: swap columns AK (37) and AL (38) --------------
$ws = $wb.Worksheets.Item("Elements")
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $akCell = $ws.Cells.Item($r, 37)
    $alCell = $ws.Cells.Item($r, 38)
    $akVal = $akCell.Value2
    $alVal = $alCell.Value2
    if ($akVal -ne $alVal) {
        $akCell.Value2 = $alVal
        $alCell.Value2 = $akVal
    }
}

# Column widths follow the content: AK becomes the wide
# "Spécification métier" column, AL the narrow "RIM Mapping" one.
# (Values picked so the engine's pixel-quantised ColumnWidth lands as
# close as possible to the original 86.46875 / 24.98046875 widths.)
$ws.Columns.Item(37).ColumnWidth = 85.66666666666667
$ws.Columns.Item(38).ColumnWidth = 24.166666666666668
